$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.461.78'
$ws.Range("E2").Value = '  +0.35%  '

# Row 3
$ws.Range("D3").Value = '3.691.83'
$ws.Range("E3").Value = '  +0.42%  '

# Row 4
$ws.Range("E4").Value = '  +0.20%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '686.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.73%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.78'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.42%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.21%  '

# Row 8
$ws.Range("E8").Value = '  +3.25%  '

# Row 9
$ws.Range("E9").Value = '  -0.17%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.15'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.72%  '

# Row 11
$ws.Range("E11").Value = '  -2.04%  '

# Row 12
$ws.Range("E12").Value = '  +2.50%  '

# Row 13
$ws.Range("D13").Value = '4.316.26'
$ws.Range("E13").Value = '  +0.78%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.54'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.71%  '

# Row 15
$ws.Range("D15").Value = '3.698.55'
$ws.Range("E15").Value = '  +0.31%  '

# Row 16
$ws.Range("D16").Value = '69.455.60'
$ws.Range("E16").Value = '  +0.70%  '

# Row 17
$ws.Range("E17").Value = '  +2.04%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '15.90'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.50%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.46'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.40%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '471.62'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.41%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.06'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.96%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.651'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.09%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '79.81'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.57%  '

# Row 24
$ws.Range("D24").Value = '3.837.81'
$ws.Range("E24").Value = '  +0.90%  '

# Row 25
$ws.Range("E25").Value = '  -0.01%  '

# Row 26
$ws.Range("E26").Value = '  -0.38%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.03'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.05%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.28'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.11%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.73'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.17%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.74'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.59%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.01'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.43%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.62'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.35%  '

# Row 33
$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.03%  '

# Row 34
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.96'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.55%  '

# Row 35
$ws.Range("D35").Value = '3.665.54'
$ws.Range("E35").Value = '  +0.81%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.159'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.27%  '

# Row 37
$ws.Range("E37").Value = '  -2.23%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.19'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.95%  '

# Row 39
$ws.Range("E39").Value = '  +1.87%  '

# Row 40
$ws.Range("E40").Value = '  +0.03%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0906'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.36%  '

# Row 42
$ws.Range("E42").Value = '  +0.30%  '

# Row 43
$ws.Range("E43").Value = '  +0.13%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '165.95'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.42%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '47.48'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.01%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.13'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.09%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.74'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.40%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.000281'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.95%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.31'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.79%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '28.44'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.38%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.83'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.62%  '
